$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.444.26'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.917.06'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.08'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4819'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4076'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08236'
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.015'
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.41'
$ws.Range('D12').Value = '1.921.75'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.079'
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.246'
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.44'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06810'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D21').Value = '29.468.10'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.658'
$ws.Range('E22').Value = '  +2.33%  '
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.179'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '2.145.01'
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.668'
$ws.Range('E26').Value = '  +10.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.92'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.04'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.114'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.019'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09600'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.707'
$ws.Range('E33').Value = '  +6.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.555'
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02287'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06113'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.093'
$ws.Range('E39').Value = '  +3.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5986'
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.85'
$ws.Range('E41').Value = '  +7.00%  '
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.407'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.279'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07603'
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.45'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5597'
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '118.37'
$ws.Range('E49').Value = '  +4.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.431'
$ws.Range('E50').Value = '  +4.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.44'
$ws.Range('E51').Value = '  +0.74%  '
